# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (cloned from the "2022-Q2" sheet so it
# keeps identical formatting) right after the "总计" summary sheet, fills it
# in with the new quarter's fund data, and adds the corresponding summary
# row at the top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" detail sheet by cloning "2022-Q2" (so all
#    styles / column widths / page setup match the other quarterly sheets)
#    and dropping it in right before "2022-Q2" (i.e. right after 总计).
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2, $null)

$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

# Fill in the new quarter's fund row (same fund, new numbers). D/E/F/G are
# stored as text in this workbook, so the values are entered with a
# leading apostrophe to keep them as text instead of being parsed as
# numbers.
$wsQ3.Range("D2").Value = "'1.12"
$wsQ3.Range("E2").Value = "'90.06"
$wsQ3.Range("F2").Value = "'2.53"
$wsQ3.Range("G2").Value = "'0.0283"
$wsQ3.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2. Add the matching summary row to the "总计" sheet, as the new first
#    data row (row 2), shifting the existing rows down.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# The insert drags the header row's formatting into the new row; strip it
# back off the text/number cells (B:D) and restore column A's style from
# the row below so the new row matches the other data rows.
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.03

# Renumber the running index in column A for the rows that shifted down.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("A7").Value = 5

# ---------------------------------------------------------------------
# 3. Keep "2021-Q2" as the selected/active sheet, same as before the edit.
# ---------------------------------------------------------------------
$wsLast = $wb.Worksheets.Item("2021-Q2")
$wsLast.Activate()
